$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 value (phone number gets country code prefix)
$ws.Range("A2").Value = 919629459258

# Add new row 5 data
$ws.Range("A5").Value = 916385026448
$ws.Range("B5").Value = 9629522931

# Update selection to A5 (matches new active cell after edits)
$ws.Range("A5").Select()
